$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# New rows 22-25 for the R Assignments block
$ws.Range("A22").Value = "Dataframe Access"
$ws.Range("B22").Value = 17

$ws.Range("A23").Value = "Database Access"
$ws.Range("B23").Value = 10

$ws.Range("A24").Value = "Sum R"
$ws.Range("B24").Formula = "=SUM(B17:B23)"
$ws.Range("A24").Font.Bold = $true

$ws.Range("A25").Value = "Minimum R"
$ws.Range("B25").Formula = "=ROUND(B24*0.75,0)"
$ws.Range("A25").Font.Bold = $true
$ws.Range("B25").Font.Bold = $true

# B13 (Minimum Python result) becomes bold too, matching the new style
$ws.Range("B13").Font.Bold = $true

# Adjust the view: scroll position and active selection cell
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 5
$ws.Range("B13").Select()
